$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Update the SqlIP value (cell E2) from "127.0.0.1" to "192.168.0.24"
$ws.Range("E2").Value = "192.168.0.24"

# Move the active cell selection to H6 (matches final saved view state)
$ws.Range("H6").Select()
